$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 341.1
$ws.Range("I6").Value = 76.25
$ws.Range("K6").Value = 228.75
$ws.Range("M6").Value = -116.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 865.7705
$ws.Range("J129").Value = 889.38184
$ws.Range("L129").Value = 2668.14552
$ws.Range("N129").Value = -12668.14552

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1698.6825
$ws.Range("I138").Value = 1210.1936
$ws.Range("J138").Value = 2171.9062
$ws.Range("K138").Value = 3630.5808
$ws.Range("L138").Value = 6515.7186
$ws.Range("M138").Value = 1509.4192
$ws.Range("N138").Value = -16795.7186

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 560
$ws.Range("I5").Value = 90
$ws.Range("K5").Value = 90
$ws.Range("M5").Value = 22

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 4500
$ws.Range("I25").Value = 4500
$ws.Range("K25").Value = 4500
$ws.Range("M25").Value = -4098

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3308.2266
$ws.Range("I32").Value = 2003
$ws.Range("J32").Value = 10160.667
$ws.Range("K32").Value = 2003
$ws.Range("L32").Value = 10160.667
$ws.Range("M32").Value = -1716
$ws.Range("N32").Value = -10734.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1392.3077
$ws.Range("I45").Value = 1088.0625
$ws.Range("J45").Value = 1879.1
$ws.Range("K45").Value = 1088.0625
$ws.Range("L45").Value = 1879.1
$ws.Range("M45").Value = -711.0625
$ws.Range("N45").Value = -2633.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 560
$ws.Range("I4").Value = 90
$ws.Range("K4").Value = 90
$ws.Range("M4").Value = 25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 229.66667
$ws.Range("I22").Value = 229.66667
$ws.Range("K22").Value = 229.66667
$ws.Range("M22").Value = -56.66667000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H68").Value = 27375
$ws.Range("I68").Value = 34750
$ws.Range("K68").Value = 34750
$ws.Range("M68").Value = -33939

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H71").Value = 27375
$ws.Range("I71").Value = 34750
$ws.Range("K71").Value = 104250
$ws.Range("M71").Value = -100194

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 20562.5
$ws.Range("I82").Value = 10750
$ws.Range("K82").Value = 10750
$ws.Range("M82").Value = -10367

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value = 20562.5
$ws.Range("I85").Value = 10750
$ws.Range("K85").Value = 10750
$ws.Range("M85").Value = -9424

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 98062.86
$ws.Range("I86").Value = 3776.5
$ws.Range("J86").Value = 183777.73
$ws.Range("K86").Value = 3776.5
$ws.Range("L86").Value = 183777.73
$ws.Range("M86").Value = -2653.5
$ws.Range("N86").Value = -186023.73

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H87").Value = 40000
$ws.Range("I87").Value = 40000
$ws.Range("K87").Value = 40000
$ws.Range("M87").Value = -38752

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 98062.86
$ws.Range("I89").Value = 3776.5
$ws.Range("J89").Value = 183777.73
$ws.Range("K89").Value = 18882.5
$ws.Range("L89").Value = 918888.65
$ws.Range("M89").Value = -13266.5
$ws.Range("N89").Value = -930120.65

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H90").Value = 40000
$ws.Range("I90").Value = 40000
$ws.Range("K90").Value = 120000
$ws.Range("M90").Value = -113760

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 695.8125
$ws.Range("I94").Value = 741
$ws.Range("K94").Value = 741
$ws.Range("M94").Value = -290

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1132.1818
$ws.Range("I107").Value = 1183.6666
$ws.Range("K107").Value = 1183.6666
$ws.Range("M107").Value = 736.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2566.175
$ws.Range("I134").Value = 2391.6333
$ws.Range("J134").Value = 3089.8
$ws.Range("K134").Value = 7174.8999
$ws.Range("L134").Value = 9269.400000000001
$ws.Range("M134").Value = -4639.8999
$ws.Range("N134").Value = -14339.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1036.3462
$ws.Range("I134").Value = 1008.7727
$ws.Range("K134").Value = 3026.3181
$ws.Range("M134").Value = -491.3181

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 359.57144
$ws.Range("J7").Value = 528.8570999999999
$ws.Range("L7").Value = 1586.5713
$ws.Range("N7").Value = -1810.5713

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 1750
$ws.Range("J80").Value = 1750
$ws.Range("L80").Value = 5250
$ws.Range("N80").Value = -7122

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 1750
$ws.Range("J83").Value = 1750
$ws.Range("L83").Value = 15750
$ws.Range("N83").Value = -25110

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 727.7778
$ws.Range("J117").Value = 954.6
$ws.Range("L117").Value = 2863.8
$ws.Range("N117").Value = -9747.799999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 797.0599999999999
$ws.Range("J131").Value = 806.0947
$ws.Range("L131").Value = 2418.2841
$ws.Range("N131").Value = -12498.2841

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 6500
$ws.Range("I99").Value = 6500
$ws.Range("K99").Value = 6500
$ws.Range("M99").Value = -4254

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2178.4783
$ws.Range("I102").Value = 2243.1428
$ws.Range("K102").Value = 2243.1428
$ws.Range("M102").Value = -621.1428000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 27782784
$ws.Range("I126").Value = 27782784
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 83348352
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -83345882
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H127").Value = 33013.625
$ws.Range("J127").Value = 33013.625
$ws.Range("L127").Value = 33013.625
$ws.Range("N127").Value = -42933.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4010.5
$ws.Range("I136").Value = 3093.077
$ws.Range("K136").Value = 9279.231
$ws.Range("M136").Value = -6729.231

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 858
$ws.Range("I100").Value = 449.33334
$ws.Range("J100").Value = 1266.6666
$ws.Range("K100").Value = 898.66668
$ws.Range("L100").Value = 2533.3332
$ws.Range("M100").Value = -357.66668
$ws.Range("N100").Value = -3615.3332

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1660.4
$ws.Range("I113").Value = 1100.6666
$ws.Range("J113").Value = 2500
$ws.Range("K113").Value = 3301.9998
$ws.Range("L113").Value = 7500
$ws.Range("M113").Value = -1131.9998
$ws.Range("N113").Value = -11840

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 19843448
$ws.Range("I136").Value = 27779706
$ws.Range("J136").Value = 2799.625
$ws.Range("K136").Value = 83339118
$ws.Range("L136").Value = 8398.875
$ws.Range("M136").Value = -83336568
$ws.Range("N136").Value = -13498.875
